$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.823.96"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "1.661.77"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'215.36"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  +5.08%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").Value = "'0.0619"
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "'20.16"
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").Value = "'0.0895"
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("D12").Value = "1.894.76"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "1.659.67"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "'65.98"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "26.826.26"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "'232.05"
$ws.Range("E18").Value = "  -2.83%  "
$ws.Range("D19").Value = "'7.86"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "0.0₃0730"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'4.42"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("D24").Value = "'9.16"
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").Value = "'145.73"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").Value = "'15.86"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "'3.34"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").Value = "1.458.84"
$ws.Range("E33").Value = "  -5.19%  "
$ws.Range("D34").Value = "'3.15"
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'0.573"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'0.897"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "'5.83"
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("E43").Value = "  +6.23%  "
$ws.Range("D44").Value = "'65.71"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "1.807.75"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'0.777"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").Value = "'90.34"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("E51").Value = "  +0.38%  "
